$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 179
$ws.Cells.Item(2, 9).Value = 173.93333
$ws.Cells.Item(2, 10).Value = 217
$ws.Cells.Item(2, 11).Value = 173.93333
$ws.Cells.Item(2, 12).Value = 217
$ws.Cells.Item(2, 13).Value = -60.93333000000001
$ws.Cells.Item(2, 14).Value = -443

$ws.Cells.Item(11, 8).Value = 74.8
$ws.Cells.Item(11, 9).Value = 74.8
$ws.Cells.Item(11, 11).Value = 74.8
$ws.Cells.Item(11, 13).Value = 65.2

$ws.Cells.Item(64, 8).Value = 4057.8333
$ws.Cells.Item(64, 10).Value = 3889.75
$ws.Cells.Item(64, 12).Value = 3889.75
$ws.Cells.Item(64, 14).Value = -4385.75

$ws.Cells.Item(67, 8).Value = 4057.8333
$ws.Cells.Item(67, 10).Value = 3889.75
$ws.Cells.Item(67, 12).Value = 3889.75
$ws.Cells.Item(67, 14).Value = -5605.75

$ws.Cells.Item(76, 8).Value = 6250
$ws.Cells.Item(76, 9).Value = 6666.6665
$ws.Cells.Item(76, 10).Value = 5000
$ws.Cells.Item(76, 11).Value = 6666.6665
$ws.Cells.Item(76, 12).Value = 5000
$ws.Cells.Item(76, 13).Value = -6351.6665
$ws.Cells.Item(76, 14).Value = -5630

$ws.Cells.Item(79, 8).Value = 6250
$ws.Cells.Item(79, 9).Value = 6666.6665
$ws.Cells.Item(79, 10).Value = 5000
$ws.Cells.Item(79, 11).Value = 6666.6665
$ws.Cells.Item(79, 12).Value = 5000
$ws.Cells.Item(79, 13).Value = -5574.6665
$ws.Cells.Item(79, 14).Value = -7184

$ws.Cells.Item(86, 8).Value = 1985.0714
$ws.Cells.Item(86, 9).Value = 1959.1
$ws.Cells.Item(86, 10).Value = 2050
$ws.Cells.Item(86, 11).Value = 1959.1
$ws.Cells.Item(86, 12).Value = 2050
$ws.Cells.Item(86, 13).Value = -836.0999999999999
$ws.Cells.Item(86, 14).Value = -4296

$ws.Cells.Item(89, 8).Value = 1985.0714
$ws.Cells.Item(89, 9).Value = 1959.1
$ws.Cells.Item(89, 10).Value = 2050
$ws.Cells.Item(89, 11).Value = 9795.5
$ws.Cells.Item(89, 12).Value = 10250
$ws.Cells.Item(89, 13).Value = -4179.5
$ws.Cells.Item(89, 14).Value = -21482

$ws.Cells.Item(107, 8).Value = 2320.72
$ws.Cells.Item(107, 9).Value = 1862.2778
$ws.Cells.Item(107, 10).Value = 3499.5715
$ws.Cells.Item(107, 11).Value = 1862.2778
$ws.Cells.Item(107, 12).Value = 3499.5715
$ws.Cells.Item(107, 13).Value = 57.72219999999993
$ws.Cells.Item(107, 14).Value = -7339.5715

$ws.Cells.Item(112, 8).Value = 2568
$ws.Cells.Item(112, 10).Value = 3066.3333
$ws.Cells.Item(112, 12).Value = 9198.999899999999
$ws.Cells.Item(112, 14).Value = -11414.9999

$ws.Cells.Item(125, 8).Value = 1781.5714
$ws.Cells.Item(125, 10).Value = 2880.5454
$ws.Cells.Item(125, 12).Value = 25924.9086
$ws.Cells.Item(125, 14).Value = -30844.9086

$ws.Cells.Item(138, 8).Value = 469284.66
$ws.Cells.Item(138, 10).Value = 597428.2
$ws.Cells.Item(138, 12).Value = 1792284.6
$ws.Cells.Item(138, 14).Value = -1802564.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6467.9424
$ws.Cells.Item(32, 9).Value = 5359.174
$ws.Cells.Item(32, 10).Value = 10718.223
$ws.Cells.Item(32, 11).Value = 5359.174
$ws.Cells.Item(32, 12).Value = 10718.223
$ws.Cells.Item(32, 13).Value = -5072.174
$ws.Cells.Item(32, 14).Value = -11292.223

$ws.Cells.Item(63, 8).Value = 30305182
$ws.Cells.Item(63, 9).Value = 2110.7368
$ws.Cells.Item(63, 10).Value = 71430776
$ws.Cells.Item(63, 11).Value = 2110.7368
$ws.Cells.Item(63, 12).Value = 71430776
$ws.Cells.Item(63, 13).Value = -1424.7368
$ws.Cells.Item(63, 14).Value = -71432148

$ws.Cells.Item(66, 8).Value = 30305182
$ws.Cells.Item(66, 9).Value = 2110.7368
$ws.Cells.Item(66, 10).Value = 71430776
$ws.Cells.Item(66, 11).Value = 10553.684
$ws.Cells.Item(66, 12).Value = 357153880
$ws.Cells.Item(66, 13).Value = -7121.684000000001
$ws.Cells.Item(66, 14).Value = -357160744

$ws.Cells.Item(88, 8).Value = 2661.4
$ws.Cells.Item(88, 10).Value = 2826.75
$ws.Cells.Item(88, 12).Value = 2826.75
$ws.Cells.Item(88, 14).Value = -3638.75

$ws.Cells.Item(91, 8).Value = 2661.4
$ws.Cells.Item(91, 10).Value = 2826.75
$ws.Cells.Item(91, 12).Value = 2826.75
$ws.Cells.Item(91, 14).Value = -5634.75

$ws.Cells.Item(122, 8).Value = 1800.1818
$ws.Cells.Item(122, 9).Value = 1522.1765
$ws.Cells.Item(122, 11).Value = 4566.529500000001
$ws.Cells.Item(122, 13).Value = -2116.529500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 408.45456
$ws.Cells.Item(64, 10).Value = 399.25
$ws.Cells.Item(64, 12).Value = 399.25
$ws.Cells.Item(64, 14).Value = -849.25

$ws.Cells.Item(67, 8).Value = 408.45456
$ws.Cells.Item(67, 10).Value = 399.25
$ws.Cells.Item(67, 12).Value = 399.25
$ws.Cells.Item(67, 14).Value = -1959.25

$ws.Cells.Item(105, 8).Value = 166668350
$ws.Cells.Item(105, 9).Value = 166668350
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 166668350
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).ClearContents()
$ws.Cells.Item(105, 14).Value = -166666603

$ws.Cells.Item(134, 8).Value = 1478.1
$ws.Cells.Item(134, 9).Value = 1309
$ws.Cells.Item(134, 11).Value = 3927
$ws.Cells.Item(134, 13).Value = -1392

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 10002220
$ws.Cells.Item(62, 9).Value = 2341.1765
$ws.Cells.Item(62, 11).Value = 2341.1765
$ws.Cells.Item(62, 13).Value = -1717.1765

$ws.Cells.Item(65, 8).Value = 10002220
$ws.Cells.Item(65, 9).Value = 2341.1765
$ws.Cells.Item(65, 11).Value = 11705.8825
$ws.Cells.Item(65, 13).Value = -8585.8825

$ws.Cells.Item(99, 8).Value = 1599.1
$ws.Cells.Item(99, 9).Value = 1576.7778
$ws.Cells.Item(99, 10).Value = 1800
$ws.Cells.Item(99, 11).Value = 1576.7778
$ws.Cells.Item(99, 12).Value = 1800
$ws.Cells.Item(99, 13).Value = -78.77780000000007
$ws.Cells.Item(99, 14).Value = -4796

$ws.Cells.Item(126, 8).Value = 1599.1
$ws.Cells.Item(126, 9).Value = 1576.7778
$ws.Cells.Item(126, 10).Value = 1800
$ws.Cells.Item(126, 11).Value = 4730.3334
$ws.Cells.Item(126, 12).Value = 5400
$ws.Cells.Item(126, 13).Value = -2260.3334
$ws.Cells.Item(126, 14).Value = -10340

$ws.Cells.Item(132, 8).Value = 2486.35
$ws.Cells.Item(132, 9).Value = 2048.6
$ws.Cells.Item(132, 11).Value = 6145.799999999999
$ws.Cells.Item(132, 13).Value = -3615.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 638.8823
$ws.Cells.Item(2, 9).Value = 21.166666
$ws.Cells.Item(2, 10).Value = 2121.4
$ws.Cells.Item(2, 11).Value = 126.999996
$ws.Cells.Item(2, 12).Value = 12728.4
$ws.Cells.Item(2, 13).Value = -13.999996
$ws.Cells.Item(2, 14).Value = -12954.4

$ws.Cells.Item(80, 8).Value = 3386.9
$ws.Cells.Item(80, 10).Value = 3249.3684
$ws.Cells.Item(80, 12).Value = 9748.1052
$ws.Cells.Item(80, 14).Value = -11620.1052

$ws.Cells.Item(83, 8).Value = 3386.9
$ws.Cells.Item(83, 10).Value = 3249.3684
$ws.Cells.Item(83, 12).Value = 29244.3156
$ws.Cells.Item(83, 14).Value = -38604.3156

$ws.Cells.Item(122, 8).Value = 1672.65
$ws.Cells.Item(122, 10).Value = 1926.9333
$ws.Cells.Item(122, 12).Value = 17342.3997
$ws.Cells.Item(122, 14).Value = -22242.3997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 64288780
$ws.Cells.Item(70, 9).Value = 62503388
$ws.Cells.Item(70, 10).Value = 66669304
$ws.Cells.Item(70, 11).Value = 62503388
$ws.Cells.Item(70, 12).Value = 66669304
$ws.Cells.Item(70, 13).Value = -62503118
$ws.Cells.Item(70, 14).Value = -66669844

$ws.Cells.Item(73, 8).Value = 64288780
$ws.Cells.Item(73, 9).Value = 62503388
$ws.Cells.Item(73, 10).Value = 66669304
$ws.Cells.Item(73, 11).Value = 62503388
$ws.Cells.Item(73, 12).Value = 66669304
$ws.Cells.Item(73, 13).Value = -62502452
$ws.Cells.Item(73, 14).Value = -66671176

$ws.Cells.Item(80, 8).Value = 9480
$ws.Cells.Item(80, 9).Value = 20000
$ws.Cells.Item(80, 10).Value = 6850
$ws.Cells.Item(80, 11).Value = 20000
$ws.Cells.Item(80, 12).Value = 6850
$ws.Cells.Item(80, 13).Value = -19002
$ws.Cells.Item(80, 14).Value = -8846

$ws.Cells.Item(83, 8).Value = 9480
$ws.Cells.Item(83, 9).Value = 20000
$ws.Cells.Item(83, 10).Value = 6850
$ws.Cells.Item(83, 11).Value = 100000
$ws.Cells.Item(83, 12).Value = 34250
$ws.Cells.Item(83, 13).Value = -95008
$ws.Cells.Item(83, 14).Value = -44234

$ws.Cells.Item(126, 8).Value = 1837.25
$ws.Cells.Item(126, 9).Value = 1605.1
$ws.Cells.Item(126, 10).Value = 2069.4
$ws.Cells.Item(126, 11).Value = 4815.299999999999
$ws.Cells.Item(126, 12).Value = 6208.200000000001
$ws.Cells.Item(126, 13).Value = -2345.299999999999
$ws.Cells.Item(126, 14).Value = -11148.2

$ws.Cells.Item(132, 8).Value = 2697.303
$ws.Cells.Item(132, 9).Value = 2391.96
$ws.Cells.Item(132, 11).Value = 7175.88
$ws.Cells.Item(132, 13).Value = -4645.88

$ws.Cells.Item(135, 8).Value = 34385.863
$ws.Cells.Item(135, 10).Value = 33828.215
$ws.Cells.Item(135, 12).Value = 33828.215
$ws.Cells.Item(135, 14).Value = -43968.215

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2797
$ws.Cells.Item(7, 9).Value = 2594.25
$ws.Cells.Item(7, 10).Value = 3202.5
$ws.Cells.Item(7, 11).Value = 2594.25
$ws.Cells.Item(7, 12).Value = 3202.5
$ws.Cells.Item(7, 13).Value = -2482.25
$ws.Cells.Item(7, 14).Value = -3426.5

$ws.Cells.Item(40, 8).Value = 4866.25
$ws.Cells.Item(40, 9).Value = 2332.6667
$ws.Cells.Item(40, 10).Value = 6386.4
$ws.Cells.Item(40, 11).Value = 2332.6667
$ws.Cells.Item(40, 12).Value = 6386.4
$ws.Cells.Item(40, 13).Value = -2196.6667
$ws.Cells.Item(40, 14).Value = -6658.4

$ws.Cells.Item(43, 8).Value = 4000
$ws.Cells.Item(43, 10).Value = 4000
$ws.Cells.Item(43, 12).Value = 4000
$ws.Cells.Item(43, 14).Value = -4386

$ws.Cells.Item(122, 8).Value = 17858804
$ws.Cells.Item(122, 9).Value = 22728622
$ws.Cells.Item(122, 10).Value = 2801.6667
$ws.Cells.Item(122, 11).Value = 68185866
$ws.Cells.Item(122, 12).Value = 8405.000100000001
$ws.Cells.Item(122, 13).Value = -68183416
$ws.Cells.Item(122, 14).Value = -13305.0001

$ws.Cells.Item(126, 8).Value = 2797
$ws.Cells.Item(126, 9).Value = 2594.25
$ws.Cells.Item(126, 10).Value = 3202.5
$ws.Cells.Item(126, 11).Value = 7782.75
$ws.Cells.Item(126, 12).Value = 9607.5
$ws.Cells.Item(126, 13).Value = -5312.75
$ws.Cells.Item(126, 14).Value = -14547.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 13).ClearContents()

$ws.Cells.Item(122, 8).Value = 10418456
$ws.Cells.Item(122, 9).Value = 11365474
$ws.Cells.Item(122, 11).Value = 34096422
$ws.Cells.Item(122, 13).Value = -34093972

$ws.Cells.Item(136, 8).Value = 1435.75
$ws.Cells.Item(136, 9).Value = 1126.4166
$ws.Cells.Item(136, 11).Value = 3379.2498
$ws.Cells.Item(136, 13).Value = -829.2498000000001
